$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in progress/status notes in column B for select rows
$ws.Range("B3").Value = "in progress"
$ws.Range("B4").Value = "done"
$ws.Range("B9").Value = "done"
$ws.Range("B10").Value = "20 enemies " + [char]0x2026 + " porbably done"
$ws.Range("B11").Value = "possibly done I just changed the colors of the stars we could do different images though"
$ws.Range("B12").Value = "same as MP3 "
$ws.Range("B14").Value = "150 enemies - done"

# Widen column B to fit the new content
$ws.Columns.Item(2).ColumnWidth = 47.83

# Move the active selection to A22
$ws.Range("A22").Select()
